# Refresh Leve profit calculations (currentAveragePrice / LevePrice / LeveProfit
# columns H-N) across all job sheets, sourced from updated Universalis market
# data during the scheduled data-refresh run.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3085.3684
$ws.Range("J32").Value = 4515.625
$ws.Range("L32").Value = 4515.625
$ws.Range("N32").Value = -5167.625
$ws.Range("H62").Value = 5427.9165
$ws.Range("I62").Value = 5244.5713
$ws.Range("K62").Value = 5244.5713
$ws.Range("M62").Value = -4620.5713
$ws.Range("H65").Value = 5427.9165
$ws.Range("I65").Value = 5244.5713
$ws.Range("K65").Value = 26222.8565
$ws.Range("M65").Value = -23102.8565
$ws.Range("H138").Value = 4833.0884
$ws.Range("J138").Value = 3556.4285
$ws.Range("L138").Value = 10669.2855
$ws.Range("N138").Value = -20949.2855

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2586.75
$ws.Range("I32").Value = 972.7027
$ws.Range("K32").Value = 972.7027
$ws.Range("M32").Value = -685.7027
$ws.Range("H45").Value = 39256.258
$ws.Range("I45").Value = 47951.137
$ws.Range("K45").Value = 47951.137
$ws.Range("M45").Value = -47574.137
$ws.Range("H61").Value = 1276221
$ws.Range("I61").Value = 29149.316
$ws.Range("J61").Value = 6541635
$ws.Range("K61").Value = 29149.316
$ws.Range("L61").Value = 6541635
$ws.Range("M61").Value = -28937.316
$ws.Range("N61").Value = -6542059
$ws.Range("H74").Value = 560361.2
$ws.Range("I74").Value = 2625.2273
$ws.Range("K74").Value = 2625.2273
$ws.Range("M74").Value = -1751.2273
$ws.Range("H77").Value = 560361.2
$ws.Range("I77").Value = 2625.2273
$ws.Range("K77").Value = 13126.1365
$ws.Range("M77").Value = -8758.136500000001
$ws.Range("H88").Value = 1277.3334
$ws.Range("I88").Value = 1576.6666
$ws.Range("J88").Value = 1127.6666
$ws.Range("K88").Value = 1576.6666
$ws.Range("L88").Value = 1127.6666
$ws.Range("M88").Value = -1170.6666
$ws.Range("N88").Value = -1939.6666
$ws.Range("H91").Value = 1277.3334
$ws.Range("I91").Value = 1576.6666
$ws.Range("J91").Value = 1127.6666
$ws.Range("K91").Value = 1576.6666
$ws.Range("L91").Value = 1127.6666
$ws.Range("M91").Value = -172.6666
$ws.Range("N91").Value = -3935.6666
$ws.Range("H97").Value = 5310.3477
$ws.Range("I97").Value = 5415.143
$ws.Range("K97").Value = 5415.143
$ws.Range("M97").Value = -4919.143
$ws.Range("H122").Value = 2520.5518
$ws.Range("I122").Value = 2318.1304
$ws.Range("K122").Value = 6954.3912
$ws.Range("M122").Value = -4504.3912
$ws.Range("H136").Value = 1276221
$ws.Range("I136").Value = 29149.316
$ws.Range("J136").Value = 6541635
$ws.Range("K136").Value = 87447.948
$ws.Range("L136").Value = 19624905
$ws.Range("M136").Value = -84897.948
$ws.Range("N136").Value = -19630005

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4803
$ws.Range("I86").Value = 2524.4
$ws.Range("K86").Value = 2524.4
$ws.Range("M86").Value = -1401.4
$ws.Range("H89").Value = 4803
$ws.Range("I89").Value = 2524.4
$ws.Range("K89").Value = 12622
$ws.Range("M89").Value = -7006
$ws.Range("H107").Value = 10585.885
$ws.Range("I107").Value = 13270.263
$ws.Range("K107").Value = 13270.263
$ws.Range("M107").Value = -11350.263
$ws.Range("H134").Value = 33336234
$ws.Range("I134").Value = 2678.5264
$ws.Range("K134").Value = 8035.5792
$ws.Range("M134").Value = -5500.5792

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1146.1666
$ws.Range("J107").Value = 1241.2727
$ws.Range("L107").Value = 1241.2727
$ws.Range("N107").Value = -5081.2727
$ws.Range("H134").Value = 2778.6316
$ws.Range("I134").Value = 2216.5386
$ws.Range("K134").Value = 6649.6158
$ws.Range("M134").Value = -4114.6158

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 5345.769
$ws.Range("I22").Value = 7212
$ws.Range("K22").Value = 21636
$ws.Range("M22").Value = -21467
$ws.Range("H27").Value = 5345.769
$ws.Range("I27").Value = 7212
$ws.Range("K27").Value = 21636
$ws.Range("M27").Value = -21534
$ws.Range("H42").Value = 18000
$ws.Range("J42").Value = 18000
$ws.Range("L42").Value = 54000
$ws.Range("N42").Value = -55068
$ws.Range("H44").Value = 4399.4546
$ws.Range("I44").Value = 1484.8572
$ws.Range("K44").Value = 4454.571599999999
$ws.Range("M44").Value = -4056.571599999999
$ws.Range("H50").Value = 1598
$ws.Range("I50").Value = 1953.1428
$ws.Range("K50").Value = 5859.428400000001
$ws.Range("M50").Value = -5378.428400000001
$ws.Range("H53").Value = 1598
$ws.Range("I53").Value = 1953.1428
$ws.Range("K53").Value = 5859.428400000001
$ws.Range("M53").Value = -5378.428400000001
$ws.Range("H56").Value = 10996410
$ws.Range("I56").Value = 10996410
$ws.Range("K56").Value = 10996410
$ws.Range("M56").Value = -10995880
$ws.Range("H58").Value = 9852.941000000001
$ws.Range("J58").Value = 10312.5
$ws.Range("L58").Value = 30937.5
$ws.Range("N58").Value = -31193.5
$ws.Range("H94").Value = 932927.4399999999
$ws.Range("I94").Value = 2024
$ws.Range("K94").Value = 6072
$ws.Range("M94").Value = -5396
$ws.Range("H109").Value = 3707.9
$ws.Range("I109").Value = 786.55554
$ws.Range("K109").Value = 2359.66662
$ws.Range("M109").Value = -1319.66662
$ws.Range("H131").Value = 14432604
$ws.Range("I131").Value = 11365308
$ws.Range("K131").Value = 34095924
$ws.Range("M131").Value = -34090884
$ws.Range("H140").Value = 3905
$ws.Range("I140").Value = 2999.7
$ws.Range("K140").Value = 8999.099999999999
$ws.Range("M140").Value = -3819.099999999999

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 33998.332
$ws.Range("J46").Value = 33998.332
$ws.Range("L46").Value = 33998.332
$ws.Range("N46").Value = -34310.332
$ws.Range("H97").Value = 1570
$ws.Range("I97").Value = 1215.75
$ws.Range("J97").Value = 2987
$ws.Range("K97").Value = 1215.75
$ws.Range("L97").Value = 2987
$ws.Range("M97").Value = -719.75
$ws.Range("N97").Value = -3979
$ws.Range("H102").Value = 41670350
$ws.Range("I102").Value = 50003620
$ws.Range("K102").Value = 50003620
$ws.Range("M102").Value = -50001998
$ws.Range("H113").Value = 4419.2
$ws.Range("J113").Value = 4998.5
$ws.Range("L113").Value = 4998.5
$ws.Range("N113").Value = -9338.5
$ws.Range("H122").Value = 5849.7
$ws.Range("J122").Value = 2499.3333
$ws.Range("L122").Value = 7497.999899999999
$ws.Range("N122").Value = -12397.9999
$ws.Range("H126").Value = 3596.1428
$ws.Range("I126").Value = 3362.1667
$ws.Range("K126").Value = 10086.5001
$ws.Range("M126").Value = -7616.500100000001
$ws.Range("H135").Value = 169990
$ws.Range("J135").Value = 169990
$ws.Range("L135").Value = 169990
$ws.Range("N135").Value = -180130

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 9375
$ws.Range("I34").Value = 6250
$ws.Range("J34").Value = 12500
$ws.Range("K34").Value = 6250
$ws.Range("L34").Value = 12500
$ws.Range("M34").Value = -6078
$ws.Range("N34").Value = -12844
$ws.Range("H122").Value = 2935.919
$ws.Range("I122").Value = 2768
$ws.Range("K122").Value = 8304
$ws.Range("M122").Value = -5854
$ws.Range("H136").Value = 2615.36
$ws.Range("I136").Value = 2241.5881
$ws.Range("J136").Value = 3409.625
$ws.Range("K136").Value = 6724.7643
$ws.Range("L136").Value = 10228.875
$ws.Range("M136").Value = -4174.7643
$ws.Range("N136").Value = -15328.875

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 58479.777
$ws.Range("I81").Value = 2175.9333
$ws.Range("J81").Value = 339999
$ws.Range("K81").Value = 4351.8666
$ws.Range("L81").Value = 679998
$ws.Range("M81").Value = -3290.8666
$ws.Range("N81").Value = -682120
$ws.Range("H84").Value = 58479.777
$ws.Range("I84").Value = 2175.9333
$ws.Range("J84").Value = 339999
$ws.Range("K84").Value = 21759.333
$ws.Range("L84").Value = 3399990
$ws.Range("M84").Value = -16455.333
$ws.Range("N84").Value = -3410598
$ws.Range("H96").Value = 2429.25
$ws.Range("I96").Value = 1990
$ws.Range("J96").Value = 2492
$ws.Range("K96").Value = 1990
$ws.Range("L96").Value = 2492
$ws.Range("M96").Value = -617
$ws.Range("N96").Value = -5238
$ws.Range("H126").Value = 8336446
$ws.Range("I126").Value = 12502891
$ws.Range("K126").Value = 37508673
$ws.Range("M126").Value = -37506203
$ws.Range("H132").Value = 2094.261
$ws.Range("I132").Value = 1195.7858
$ws.Range("K132").Value = 3587.3574
$ws.Range("M132").Value = -1057.3574
